# Dashboard Calculations Figures - date revisions
# 1) Bump the cached "datetimeFigureOut" footer field text on the slide
#    master and every slide layout from the 2-digit-year form (12/1/25)
#    to the 4-digit-year form (12/1/2025).
# 2) Update the waitlisting callout on slide 1: the "last dialysis start"
#    date moves from 12/31/2023 to 3/15/2022.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $phType = $null
        try { $phType = $sh.PlaceholderFormat.Type } catch { $phType = $null }
        if ($phType -eq 16) {
            # ppPlaceholderDate
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "12/1/25") {
                $tr.Text = "12/1/2025"
            }
        }
    }
}

# --- Slide master ---
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# --- Every slide layout on the master ---
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# --- Slide 1: waitlisting callout text ---
$s = $p.Slides.Item(1)
$topShape = $s.Shapes.Item(1)
$groupItems = $topShape.GroupItems

for ($i = 1; $i -le $groupItems.Count; $i++) {
    $sh = $groupItems.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        $full = $tr.Text
        if ($full -like "*Waitlisting through 3/15/2024, last dialysis start 12/31/2023*") {
            $paras = $tr
            $found = $false
            $paraCount = $tr.Paragraphs().Count
            for ($pi = 1; $pi -le $paraCount; $pi++) {
                $para = $tr.Paragraphs($pi, 1)
                if ($para.Text -like "*Waitlisting through 3/15/2024, last dialysis start 12/31/2023*") {
                    $newText = $para.Text.Replace("last dialysis start 12/31/2023", "last dialysis start 3/15/2022")
                    $para.Text = $newText
                    $found = $true
                }
            }
        }
    }
}
